$d = $word.ActiveDocument

# --- Step 1: add "Bag of words" as a new run before the first drawing ---
# Directly inserting text at the drawing's anchor position would replace
# the drawing, so instead we insert a brand new (empty) paragraph before
# paragraph 1, put the text there, then merge it back into paragraph 1 by
# deleting the intervening paragraph mark. This leaves the drawing run
# intact and simply prepends a text run ahead of it in the same paragraph.
$p1 = $d.Paragraphs.Item(1)
$insPt = $d.Range($p1.Range.Start, $p1.Range.Start)
$insPt.InsertParagraphBefore()

$newPara = $d.Paragraphs.Item(1)
$textPt = $d.Range($newPara.Range.Start, $newPara.Range.Start)
$textPt.InsertBefore("Bag of words")

$mergedPara = $d.Paragraphs.Item(1)
$paraMark = $d.Range($mergedPara.Range.End - 1, $mergedPara.Range.End)
$paraMark.Delete()

# --- Step 2: mark both drawings as NoProof (adds <w:rPr><w:noProof/></w:rPr>) ---
for ($i = 1; $i -le $d.InlineShapes.Count; $i++) {
    $shape = $d.InlineShapes.Item($i)
    $shape.Range.NoProofing = $true
}
